$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "TextBox 2" shape on slide 1 that holds the attendance
# password placeholder ("Today's Attendance password" / "__________").
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*Attendance password*") {
            $shape = $candidate
            break
        }
    }
}

$tr = $shape.TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
$para2.Runs(1).Text = "bigo"
